$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.95399530092467
$ws.Range("C2").Value = 10.27326774597168
$ws.Range("D2").Value = 6.62
$ws.Range("E2").Value = 0.897160927841782

$ws.Range("B3").Value = 14.76794118325206
$ws.Range("C3").Value = 11.66104507446289
$ws.Range("D3").Value = 7.37
$ws.Range("E3").Value = 0.9784198060558088

$ws.Range("B4").Value = 2.726909101693961
$ws.Range("C4").Value = 2.046104192733765
$ws.Range("D4").Value = 1.33
$ws.Range("E4").Value = 0.9505789800433458

$ws.Range("B5").Value = 2.726811438621243
$ws.Range("C5").Value = 2.045972108840942
$ws.Range("D5").Value = 1.33
$ws.Range("E5").Value = 0.9505824899917691

$ws.Range("B6").Value = 12.5821658639182
$ws.Range("C6").Value = 7.802381992340088
$ws.Range("D6").Value = 7.08
$ws.Range("E6").Value = 0.8575052819786672

$ws.Range("B7").Value = 8.572147475475944
$ws.Range("C7").Value = 6.120565414428711
$ws.Range("D7").Value = 3.33
$ws.Range("E7").Value = 0.9302316218366755

$ws.Range("B8").Value = 9.165929267585415
$ws.Range("C8").Value = 6.710178375244141
$ws.Range("D8").Value = 5.1
$ws.Range("E8").Value = 0.9331685600847157

$ws.Range("B9").Value = 13.90801812407679
$ws.Range("C9").Value = 9.660391807556152
$ws.Range("D9").Value = 9.970000000000001
$ws.Range("E9").Value = 0.8511368282721031

$ws.Range("B10").Value = 60.14815083465264
$ws.Range("C10").Value = 33.17148971557617
$ws.Range("D10").Value = 28.94
$ws.Range("E10").Value = -0.5487347116377804
